$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.095739218761906
$ws.Range("C2").Value = 0.07558923879656732
$ws.Range("D2").Value = 0.07706284401608343
$ws.Range("E2").Value = 0.08671413816489881
$ws.Range("G2").Value = 2.035918829824624
$ws.Range("H2").Value = 1.711275127132865
$ws.Range("K2").Value = 0.6688731168793822
$ws.Range("L2").Value = 0.1945722440733917
$ws.Range("M2").Value = 0.2526993953038961
$ws.Range("N2").Value = 3.12931536791254

$ws.Range("B3").Value = 1.056663637891035
$ws.Range("C3").Value = 0.07186293479472283
$ws.Range("D3").Value = 0.07011787175080997
$ws.Range("E3").Value = 0.08689282353686956
$ws.Range("G3").Value = 2.02310137968297
$ws.Range("H3").Value = 1.710612454461014
$ws.Range("K3").Value = 0.6291726043023687
$ws.Range("L3").Value = 0.1921788495076839
$ws.Range("M3").Value = 0.2458665461190854
$ws.Range("N3").Value = 3.142365121340802

$ws.Range("B4").Value = 1.033301575326817
$ws.Range("C4").Value = 0.0695428988211404
$ws.Range("D4").Value = 0.06589112698657118
$ws.Range("E4").Value = 0.08702599840952807
$ws.Range("G4").Value = 2.016158212270881
$ws.Range("H4").Value = 1.710835016779072
$ws.Range("K4").Value = 0.6051830601983568
$ws.Range("L4").Value = 0.1908043155708938
$ws.Range("M4").Value = 0.241812033927733
$ws.Range("N4").Value = 3.151142766327069

$ws.Range("B5").Value = 1.023940025597938
$ws.Range("C5").Value = 0.06858935278481226
$ws.Range("D5").Value = 0.06417806874355847
$ws.Range("E5").Value = 0.08708617955949016
$ws.Range("G5").Value = 2.013561672006844
$ws.Range("H5").Value = 1.711084034860193
$ws.Range("K5").Value = 0.5955043531709521
$ws.Range("L5").Value = 0.1902681017430865
$ws.Range("M5").Value = 0.2401952474411111
$ws.Range("N5").Value = 3.154912024836065

$ws.Range("B6").Value = 1.022395135446345
$ws.Range("C6").Value = 0.06843052580824605
$ws.Range("D6").Value = 0.06389418135618996
$ws.Range("E6").Value = 0.08709652996929762
$ws.Range("G6").Value = 2.013144576884642
$ws.Range("H6").Value = 1.71113494771518
$ws.Range("K6").Value = 0.5939030830768957
$ws.Range("L6").Value = 0.1901805093379778
$ws.Range("M6").Value = 0.2399289239970344
$ws.Range("N6").Value = 3.155549519335835

$ws.Range("B7").Value = 1.033174679520243
$ws.Range("C7").Value = 0.06953007187569682
$ws.Range("D7").Value = 0.06586798619724732
$ws.Range("E7").Value = 0.08702678608306869
$ws.Range("G7").Value = 2.016122251854725
$ws.Range("H7").Value = 1.71083773403403
$ws.Range("K7").Value = 0.6050521361254084
$ws.Range("L7").Value = 0.1907969871170252
$ws.Range("M7").Value = 0.2417900857578346
$ws.Range("N7").Value = 3.151192821381159

$ws.Range("B8").Value = 1.082135162764786
$ws.Range("C8").Value = 0.07431103544578832
$ws.Range("D8").Value = 0.07466040038553956
$ws.Range("E8").Value = 0.08677088644627773
$ws.Range("G8").Value = 2.031306889474948
$ws.Range("H8").Value = 1.710916017433817
$ws.Range("K8").Value = 0.655104087280705
$ws.Range("L8").Value = 0.1937272946944262
$ws.Range("M8").Value = 0.2503142053205387
$ws.Range("N8").Value = 3.133656089177492

$ws.Range("B9").Value = 1.183149938271299
$ws.Range("C9").Value = 0.08343455523466048
$ws.Range("D9").Value = 0.09220318392999616
$ws.Range("E9").Value = 0.08645474217855487
$ws.Range("G9").Value = 2.068450366553435
$ws.Range("H9").Value = 1.716063795323691
$ws.Range("K9").Value = 0.7563313214446623
$ws.Range("L9").Value = 0.2002269493007773
$ws.Range("M9").Value = 0.2681475981416099
$ws.Range("N9").Value = 3.105340874444664

$ws.Range("B10").Value = 1.260426420810063
$ws.Range("C10").Value = 0.0899880255726373
$ws.Range("D10").Value = 0.1052812942930927
$ws.Range("E10").Value = 0.08633506938295454
$ws.Range("G10").Value = 2.100253778939958
$ws.Range("H10").Value = 1.722892998031426
$ws.Range("K10").Value = 0.83259631040562
$ws.Range("L10").Value = 0.2054614862379651
$ws.Range("M10").Value = 0.2819326337891326
$ws.Range("N10").Value = 3.088246823961711

$ws.Range("B11").Value = 1.296249116224203
$ws.Range("C11").Value = 0.09293777976556328
$ws.Range("D11").Value = 0.1112733680320588
$ws.Range("E11").Value = 0.08630495600524846
$ws.Range("G11").Value = 2.115707617002755
$ws.Range("H11").Value = 1.726662198922895
$ws.Range("K11").Value = 0.8677068932014436
$ws.Range("L11").Value = 0.2079425903959731
$ws.Range("M11").Value = 0.2883525281204129
$ws.Range("N11").Value = 3.081277143772596

$ws.Range("B12").Value = 1.309910514950332
$ws.Range("C12").Value = 0.09405031185993096
$ws.Range("D12").Value = 0.1135486423922032
$ws.Range("E12").Value = 0.08629704030947494
$ws.Range("G12").Value = 2.121701750912564
$ws.Range("H12").Value = 1.728184794927301
$ws.Range("K12").Value = 0.88106254120828
$ws.Range("L12").Value = 0.2088964723412232
$ws.Range("M12").Value = 0.2908050024959792
$ws.Range("N12").Value = 3.07875400860371

$ws.Range("B13").Value = 1.306964015044457
$ws.Range("C13").Value = 0.09381090648082591
$ws.Range("D13").Value = 0.1130583444559079
$ws.Range("E13").Value = 0.08629859015148078
$ws.Range("G13").Value = 2.120404482587929
$ws.Range("H13").Value = 1.727852638841654
$ws.Range("K13").Value = 0.8781834929611705
$ws.Range("L13").Value = 0.2086903992945111
$ws.Range("M13").Value = 0.2902758669391261
$ws.Range("N13").Value = 3.079292243344696

$ws.Range("B14").Value = 1.297371121940444
$ws.Range("C14").Value = 0.09302939788236131
$ws.Range("D14").Value = 0.1114604315285419
$ws.Range("E14").Value = 0.08630423493959682
$ws.Range("G14").Value = 2.116197908526601
$ws.Range("H14").Value = 1.726785554187757
$ws.Range("K14").Value = 0.8688044670737156
$ws.Range("L14").Value = 0.2080207795344364
$ws.Range("M14").Value = 0.288553865939015
$ws.Range("N14").Value = 3.081067235549909

$ws.Range("B15").Value = 1.291507712462987
$ws.Range("C15").Value = 0.09255011981665007
$ws.Range("D15").Value = 0.1104824735707979
$ws.Range("E15").Value = 0.08630814641165152
$ws.Range("G15").Value = 2.113639774294455
$ws.Range("H15").Value = 1.72614434270082
$ws.Range("K15").Value = 0.8630673618649496
$ws.Range("L15").Value = 0.2076124851565453
$ws.Range("M15").Value = 0.2875018766602224
$ws.Range("N15").Value = 3.082169598067665

$ws.Range("B16").Value = 1.258098780327884
$ws.Range("C16").Value = 0.08979462331421928
$ws.Range("D16").Value = 0.1048905636366868
$ws.Range("E16").Value = 0.08633752597095068
$ws.Range("G16").Value = 2.099263704367104
$ws.Range("H16").Value = 1.722660004762162
$ws.Range("K16").Value = 0.8303101501222443
$ws.Range("L16").Value = 0.2053013483399155
$ws.Range("M16").Value = 0.2815160746890086
$ws.Range("N16").Value = 3.088718551099007

$ws.Range("B17").Value = 1.237774807640847
$ws.Range("C17").Value = 0.08809619317341344
$ws.Range("D17").Value = 0.1014711029877731
$ws.Range("E17").Value = 0.0863617724375203
$ws.Range("G17").Value = 2.090697259430073
$ws.Range("H17").Value = 1.720692177590195
$ws.Range("K17").Value = 0.810321552447391
$ws.Range("L17").Value = 0.2039091084536011
$ws.Range("M17").Value = 0.2778821244912777
$ws.Range("N17").Value = 3.092942794939574

$ws.Range("B18").Value = 1.226147994305848
$ws.Range("C18").Value = 0.0871163378218256
$ws.Range("D18").Value = 0.09950834263001695
$ws.Range("E18").Value = 0.08637800870143408
$ws.Range("G18").Value = 2.085862886188579
$ws.Range("H18").Value = 1.719622698388775
$ws.Range("K18").Value = 0.7988639118628669
$ws.Range("L18").Value = 0.2031177320268682
$ws.Range("M18").Value = 0.2758059998415874
$ws.Range("N18").Value = 3.095448369007741

$ws.Range("B19").Value = 1.22222217862722
$ws.Range("C19").Value = 0.08678406526284732
$ws.Range("D19").Value = 0.09884447497131532
$ws.Range("E19").Value = 0.08638389971079441
$ws.Range("G19").Value = 2.084241984096593
$ws.Range("H19").Value = 1.719271301885954
$ws.Range("K19").Value = 0.7949912957532774
$ws.Range("L19").Value = 0.2028514012403804
$ws.Range("M19").Value = 0.275105469794596
$ws.Range("N19").Value = 3.096309744590968

$ws.Range("B20").Value = 1.239931809989685
$ws.Range("C20").Value = 0.08827730030805014
$ws.Range("D20").Value = 0.1018346936531742
$ws.Range("E20").Value = 0.08635895439499208
$ws.Range("G20").Value = 2.091599563709281
$ws.Range("H20").Value = 1.720895201780252
$ws.Range("K20").Value = 0.8124453080608873
$ws.Range("L20").Value = 0.2040563417386352
$ws.Range("M20").Value = 0.2782675131175907
$ws.Range("N20").Value = 3.092485260374062

$ws.Range("B21").Value = 1.300186179359685
$ws.Range("C21").Value = 0.09325906699849895
$ws.Range("D21").Value = 0.1119296084189898
$ws.Range("E21").Value = 0.08630248235988525
$ws.Range("G21").Value = 2.117429622538793
$ws.Range("H21").Value = 1.727096396845269
$ws.Range("K21").Value = 0.8715576848031787
$ws.Range("L21").Value = 0.2082170740344225
$ws.Range("M21").Value = 0.2890590788122722
$ws.Range("N21").Value = 3.080542724247678

$ws.Range("B22").Value = 1.340126112237556
$ws.Range("C22").Value = 0.09648887302802223
$ws.Range("D22").Value = 0.1185634181833564
$ws.Range("E22").Value = 0.08628589812550125
$ws.Range("G22").Value = 2.135139460453075
$ws.Range("H22").Value = 1.731704614183855
$ws.Range("K22").Value = 0.9105409711987988
$ws.Range("L22").Value = 0.2110199373043997
$ws.Range("M22").Value = 0.2962367123270866
$ws.Range("N22").Value = 3.073414566391079

$ws.Range("B23").Value = 1.318758194851512
$ws.Range("C23").Value = 0.09476743428710677
$ws.Range("D23").Value = 0.1150194986143447
$ws.Range("E23").Value = 0.08629289329360113
$ws.Range("G23").Value = 2.125611497179676
$ws.Range("H23").Value = 1.729194302085091
$ws.Range("K23").Value = 0.8897028276914227
$ws.Range("L23").Value = 0.2095163556248849
$ws.Range("M23").Value = 0.2923944712845525
$ws.Range("N23").Value = 3.077157001002192

$ws.Range("B24").Value = 1.238956449460773
$ws.Range("C24").Value = 0.08819543238527672
$ws.Range("D24").Value = 0.1016703045613525
$ws.Range("E24").Value = 0.08636022127760157
$ws.Range("G24").Value = 2.091191349771947
$ws.Range("H24").Value = 1.720803221867158
$ws.Range("K24").Value = 0.811485052027308
$ws.Range("L24").Value = 0.2039897494087342
$ws.Range("M24").Value = 0.2780932381696033
$ws.Range("N24").Value = 3.092691871886259

$ws.Range("B25").Value = 1.155286023230559
$ws.Range("C25").Value = 0.08099295898439607
$ws.Range("D25").Value = 0.08742456461254733
$ws.Range("E25").Value = 0.08652045156509658
$ws.Range("G25").Value = 2.0576110487454
$ws.Range("H25").Value = 1.714136268085753
$ws.Range("K25").Value = 0.7286155253617608
$ws.Range("L25").Value = 0.1983879738994361
$ws.Range("M25").Value = 0.2632033908119524
$ws.Range("N25").Value = 3.112349983198328
